$d = $word.ActiveDocument

# The text of the new paragraph, split into the same run boundaries as the
# authored edit.
$runs = @(
    "Los hilos son subprocesos dentro de un proceso que pueden llevarse ",
    "a cabo",
    " al mismo tiempo, estas tareas son ",
    "aún",
    " ",
    "más",
    " sencillas que los procesos, estos ayudan a agilizar el trabajo debido a que dentro de los procesos no se realiza una tarea a la vez sino múltiples tareas, al dividir el trabajo y realizarlo al mismo tiempo se reduce el tiempo, aumentando así la eficiencia. Si ",
    "en un ",
    "proceso no puede ",
    "realizar mas de una tarea a la vez se dice que es de único-hilo, en cambio, si un proceso realizar varias tareas a la vez se dice que es multihilo."
)

# --- 1. Insert a fresh paragraph right after the existing (only) paragraph.
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs(2)
$cur = $d.Range($newPara.Range.Start, $newPara.Range.Start)

# --- 2. Build the paragraph's runs one at a time. Inserting plain text into
# a collapsed range repeatedly would let Word coalesce same-formatted runs,
# so each run is typed into its own temporary paragraph (created with
# InsertParagraphAfter) and the separating paragraph marks are deleted
# afterwards, which leaves the runs intact but merges the paragraphs back
# into one.
$splitPositions = @()
for ($i = 0; $i -lt $runs.Count; $i++) {
    $cur.InsertAfter($runs[$i])
    if ($i -lt $runs.Count - 1) {
        $cur.InsertParagraphAfter()
        $splitPositions += $cur.End
        $cur = $d.Range($cur.End + 1, $cur.End + 1)
    }
}

# Remove the temporary paragraph marks, last-inserted first so earlier
# offsets stay valid.
for ($i = $splitPositions.Count - 1; $i -ge 0; $i--) {
    $pos = $splitPositions[$i]
    $d.Range($pos, $pos + 1).Delete()
}

# --- 3. Append the two trailing blank paragraphs. Setting Range.Text to a
# bare carriage return (rather than calling InsertParagraphAfter) produces a
# genuinely empty <w:p/> instead of one holding a placeholder empty run.
$hostPara = $d.Paragraphs(2)
$endPos = $hostPara.Range.End
$blank1 = $d.Range($endPos, $endPos)
$blank1.Text = [char]13

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endPos2 = $lastPara.Range.End
$blank2 = $d.Range($endPos2, $endPos2)
$blank2.Text = [char]13

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
